# Add 2022-Q3 sheet + data, matching commit "feat: add 2022-Q3 data"
$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1) Create the new "2022-Q3" worksheet by duplicating the "2022-Q1"
#    sheet (so it inherits identical formatting / page setup / styles)
#    and placing it right after "总计" (i.e. before "2022-Q1").
# ---------------------------------------------------------------------
$summarySheet = $wb.Worksheets.Item(1)
$templateSheet = $wb.Worksheets.Item("2022-Q1")
$templateSheet.Copy($null, $summarySheet)

$q3Sheet = $wb.Worksheets.Item(2)
$q3Sheet.Name = "2022-Q3"

# The template had 2 data rows; 2022-Q3 only needs 1, so drop row 3.
$q3Sheet.Rows.Item(3).Delete()

# Fill in the 2022-Q3 fund-holding data (row 2).
$q3Sheet.Range("B2").Value = "'090007"
$q3Sheet.Range("B2").Style = "Normal"
$q3Sheet.Range("C2").Value = "大成策略回报混合"
$q3Sheet.Range("D2").Value = "'9.90"
$q3Sheet.Range("D2").Style = "Normal"
$q3Sheet.Range("E2").Value = "'61.86"
$q3Sheet.Range("E2").Style = "Normal"
$q3Sheet.Range("F2").Value = "'2.59"
$q3Sheet.Range("F2").Style = "Normal"
$q3Sheet.Range("G2").Value = "'0.2564"
$q3Sheet.Range("G2").Style = "Normal"
$q3Sheet.Range("H2").Value = 9

# ---------------------------------------------------------------------
# 2) Update the "总计" summary sheet: insert the 2022-Q3 row at the top
#    of the data (row 2) and push the existing quarters down by one.
# ---------------------------------------------------------------------
$summarySheet.Range("A6").Copy()
$summarySheet.Range("A7").PasteSpecial(-4122)
$summarySheet.Range("A7").Value = 5
$summarySheet.Range("B7").Value = "2021-Q1"
$summarySheet.Range("C7").Value = 2
$summarySheet.Range("D7").Value = 0.78

$summarySheet.Range("B6").Value = "2021-Q2"
$summarySheet.Range("C6").Value = 4
$summarySheet.Range("D6").Value = 1.13

$summarySheet.Range("B5").Value = "2021-Q3"
$summarySheet.Range("C5").Value = 4
$summarySheet.Range("D5").Value = 1.3

$summarySheet.Range("B4").Value = "2021-Q4"
$summarySheet.Range("C4").Value = 4
$summarySheet.Range("D4").Value = 2.81

$summarySheet.Range("B3").Value = "2022-Q1"
$summarySheet.Range("C3").Value = 2
$summarySheet.Range("D3").Value = 0.24

$summarySheet.Range("B2").Value = "2022-Q3"
$summarySheet.Range("C2").Value = 1
$summarySheet.Range("D2").Value = 0.26

# ---------------------------------------------------------------------
# 3) Restore the originally-active tab ("2021-Q1", the last sheet).
# ---------------------------------------------------------------------
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$lastSheet.Select()
